# Appends the new daily COVID tracking rows (through 2021-12-08, i.e. "8/12")
# to the bottom of the existing data table on Sheet1, columns A (date serial),
# B (nuovi pos.), C (somma mobile 7gg.), D (somma mobile 7gg. per 100mila abitanti).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows to append: row, date serial, nuovi pos., somma mobile 7gg., per 100k ab.
$newData = @(
    ,@(44460, 0, 2, 37.87878787878788)
    ,@(44461, 0, 2, 37.87878787878788)
    ,@(44462, 0, 2, 37.87878787878788)
    ,@(44463, 0, 1, 18.93939393939394)
    ,@(44464, 0, 0, 0)
    ,@(44465, 0, 0, 0)
    ,@(44466, 0, 0, 0)
    ,@(44467, 0, 0, 0)
    ,@(44468, 0, 0, 0)
    ,@(44469, 0, 0, 0)
    ,@(44470, 0, 0, 0)
    ,@(44471, 0, 0, 0)
    ,@(44472, 0, 0, 0)
    ,@(44473, 0, 0, 0)
    ,@(44474, 0, 0, 0)
    ,@(44475, 0, 0, 0)
    ,@(44476, 0, 0, 0)
    ,@(44477, 0, 0, 0)
    ,@(44478, 0, 0, 0)
    ,@(44479, 0, 0, 0)
    ,@(44480, 0, 0, 0)
    ,@(44481, 0, 0, 0)
    ,@(44482, 0, 0, 0)
    ,@(44483, 0, 0, 0)
    ,@(44484, 0, 0, 0)
    ,@(44485, 0, 0, 0)
    ,@(44486, 0, 0, 0)
    ,@(44487, 0, 0, 0)
    ,@(44488, 0, 0, 0)
    ,@(44489, 0, 0, 0)
    ,@(44490, 0, 0, 0)
    ,@(44491, 0, 0, 0)
    ,@(44492, 0, 0, 0)
    ,@(44493, 0, 0, 0)
    ,@(44494, 0, 0, 0)
    ,@(44495, 0, 0, 0)
    ,@(44496, 0, 0, 0)
    ,@(44497, 0, 0, 0)
    ,@(44498, 0, 0, 0)
    ,@(44499, 0, 0, 0)
    ,@(44500, 0, 0, 0)
    ,@(44501, 0, 0, 0)
    ,@(44502, 0, 0, 0)
    ,@(44503, 0, 0, 0)
    ,@(44504, 0, 0, 0)
    ,@(44505, 0, 0, 0)
    ,@(44506, 1, 1, 18.93939393939394)
    ,@(44507, 1, 2, 37.87878787878788)
    ,@(44508, 2, 4, 75.75757575757575)
    ,@(44509, 0, 4, 75.75757575757575)
    ,@(44510, 0, 4, 75.75757575757575)
    ,@(44511, 0, 4, 75.75757575757575)
    ,@(44512, 1, 5, 94.6969696969697)
    ,@(44513, 1, 5, 94.6969696969697)
    ,@(44514, 0, 4, 75.75757575757575)
    ,@(44515, 0, 2, 37.87878787878788)
    ,@(44516, 5, 7, 132.5757575757576)
    ,@(44517, 0, 7, 132.5757575757576)
    ,@(44518, 0, 7, 132.5757575757576)
    ,@(44519, 3, 9, 170.4545454545454)
    ,@(44520, 1, 9, 170.4545454545454)
    ,@(44521, 1, 10, 189.3939393939394)
    ,@(44522, 2, 12, 227.2727272727273)
    ,@(44523, 0, 7, 132.5757575757576)
    ,@(44524, 1, 8, 151.5151515151515)
    ,@(44525, 3, 11, 208.3333333333333)
    ,@(44526, 1, 9, 170.4545454545454)
    ,@(44527, 1, 9, 170.4545454545454)
    ,@(44528, 0, 8, 151.5151515151515)
    ,@(44529, 3, 9, 170.4545454545454)
    ,@(44530, 1, 10, 189.3939393939394)
    ,@(44531, 0, 9, 170.4545454545454)
    ,@(44532, 1, 7, 132.5757575757576)
    ,@(44533, 1, 7, 132.5757575757576)
    ,@(44534, 4, 10, 189.3939393939394)
    ,@(44535, 2, 12, 227.2727272727273)
    ,@(44536, 1, 10, 189.3939393939394)
    ,@(44537, 1, 10, 189.3939393939394)
    ,@(44538, 1, 11, 208.3333333333333)
)

$firstRow = $ws.UsedRange.Rows.Count + 1
$lastRow = $firstRow + $newData.Count - 1

# Stamp the date column format (border + centered date numFmt, matching the rest of
# column A) onto the new rows by duplicating an equally sized block of already-
# formatted date cells immediately above, then fill in the real values below.
$srcFirst = $firstRow - $newData.Count
$srcLast = $firstRow - 1
$ws.Range("A" + $srcFirst + ":A" + $srcLast).Copy($ws.Range("A" + $firstRow + ":A" + $lastRow))

$r = $firstRow
foreach ($row in $newData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
